$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.918.69"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "2.007.06"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "225.37"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "0.603"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "54.39"
$ws.Range("E8").Value = "  -5.03%  "
$ws.Range("D9").Value = "0.373"
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("D10").Value = "0.0775"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").Value = "0.101"
$ws.Range("E11").Value = "  -5.46%  "
$ws.Range("D12").Value = "2.308.81"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D13").Value = "13.95"
$ws.Range("E13").Value = "  -5.56%  "
$ws.Range("D14").Value = "19.78"
$ws.Range("E14").Value = "  -5.32%  "
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "0.733"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("D17").Value = "2.060.92"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "36.837.32"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "6.33"
$ws.Range("E19").Value = "  +4.18%  "
$ws.Range("D20").Value = "68.18"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").Value = "0.0₃0808"
$ws.Range("D22").Value = "221.38"
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("E25").Value = "  -6.65%  "
$ws.Range("D26").Value = "164.86"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").Value = "9.05"
$ws.Range("E27").Value = "  -5.60%  "
$ws.Range("D28").Value = "0.125"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("D29").Value = "18.47"
$ws.Range("E29").Value = "  -2.87%  "
$ws.Range("E30").Value = "  -6.56%  "
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").Value = "0.0598"
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("D34").Value = "4.42"
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("D35").Value = "2.31"
$ws.Range("E35").Value = "  -5.63%  "
$ws.Range("D36").Value = "1.86"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  -5.45%  "
$ws.Range("D39").Value = "5.35"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "1.452.36"
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0211"
$ws.Range("E41").Value = "  -5.41%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "94.48"
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").Value = "2.77"
$ws.Range("E43").Value = "  -4.51%  "
$ws.Range("D44").Value = "0.0906"
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "15.97"
$ws.Range("E45").Value = "  -7.61%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "1.12"
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("D47").Value = "7.12"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").Value = "0.996"
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("D49").Value = "2.90"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").Value = "2.198.16"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("E51").Value = "  -9.63%  "
